$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Slow config
$ws.Cells.Item(3, 1).Value = "Slow"
$ws.Cells.Item(3, 2).Value = 5
$ws.Cells.Item(3, 3).Value = 20
$ws.Cells.Item(3, 4).Value = 3
$ws.Cells.Item(3, 5).Value = 30
$ws.Cells.Item(3, 6).Value = 3
$ws.Cells.Item(3, 7).Value = 20
$ws.Cells.Item(3, 8).Value = 2
$ws.Cells.Item(3, 9).Value = 0.8

# Row 4: Fast config
$ws.Cells.Item(4, 1).Value = "Fast"
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = 10
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 20
$ws.Cells.Item(4, 8).Value = 3
$ws.Cells.Item(4, 9).Value = 0.8

$ws.Range("I5").Select() | Out-Null
